# Auto-generated edit script for cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "1.00", "1.21") are preserved exactly as text, not coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '95.648.86'
$ws.Range("E2").Value = '  -1.51%  '
$ws.Range("D3").Value = '3.609.13'
$ws.Range("E3").Value = '  -2.21%  '
$ws.Range("E4").Value = '  +26.50%  '
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D6").Value = '222.90'
$ws.Range("E6").Value = '  -5.45%  '
$ws.Range("D7").Value = '640.20'
$ws.Range("E7").Value = '  -2.44%  '
$ws.Range("D8").Value = '0.423'
$ws.Range("E8").Value = '  -2.93%  '
$ws.Range("D9").Value = '1.21'
$ws.Range("E9").Value = '  +8.35%  '
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("D11").Value = '3.609.18'
$ws.Range("E11").Value = '  -2.19%  '
$ws.Range("D12").Value = '48.54'
$ws.Range("E12").Value = '  +8.52%  '
$ws.Range("D13").Value = '0.217'
$ws.Range("E13").Value = '  +4.60%  '
$ws.Range("E14").Value = '  -6.61%  '
$ws.Range("D15").Value = '6.52'
$ws.Range("E15").Value = '  -4.53%  '
$ws.Range("D16").Value = '4.285.15'
$ws.Range("E16").Value = '  -2.11%  '
$ws.Range("D17").Value = '24.96'
$ws.Range("E17").Value = '  +33.92%  '
$ws.Range("D18").Value = '95.672.71'
$ws.Range("E18").Value = '  -1.19%  '
$ws.Range("D19").Value = '9.03'
$ws.Range("E19").Value = '  +4.21%  '
$ws.Range("E20").Value = '  +5.94%  '
$ws.Range("D21").Value = '3.615.26'
$ws.Range("E21").Value = '  -2.24%  '
$ws.Range("D22").Value = '0.293'
$ws.Range("E22").Value = '  +42.23%  '
$ws.Range("D23").Value = '0.535'
$ws.Range("E23").Value = '  -1.53%  '
$ws.Range("D24").Value = '137.35'
$ws.Range("E24").Value = '  +22.67%  '
$ws.Range("D25").Value = '528.32'
$ws.Range("E25").Value = '  +1.52%  '
$ws.Range("E26").Value = '  -4.89%  '
$ws.Range("E27").Value = '  -8.19%  '
$ws.Range("E28").Value = '  -0.29%  '
$ws.Range("D29").Value = '3.783.63'
$ws.Range("E29").Value = '  -2.72%  '
$ws.Range("D30").Value = '12.92'
$ws.Range("E30").Value = '  -3.34%  '
$ws.Range("D31").Value = '13.28'
$ws.Range("E31").Value = '  +5.49%  '
$ws.Range("D32").Value = '3.15'
$ws.Range("E32").Value = '  +4.55%  '
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").Value = '0.636'
$ws.Range("E34").Value = '  +6.99%  '
$ws.Range("D35").Value = '33.54'
$ws.Range("E35").Value = '  +2.06%  '
$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").Value = '1.83'
$ws.Range("E36").Value = '  +1.04%  '
$ws.Range("B37").Value = 'Cronos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D37").Value = '0.183'
$ws.Range("E37").Value = '  -2.60%  '
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("B39").Value = 'USDe'
$ws.Range("C39").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '0.531'
$ws.Range("E40").Value = '  +7.40%  '
$ws.Range("D41").Value = '7.21'
$ws.Range("E41").Value = '  +5.12%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").Value = '588.30'
$ws.Range("E42").Value = '  -6.84%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").Value = '8.39'
$ws.Range("E43").Value = '  -3.75%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0532'
$ws.Range("E44").Value = '  +18.09%  '
$ws.Range("D45").Value = '41.37'
$ws.Range("E45").Value = '  +2.57%  '
$ws.Range("D46").Value = '1.01'
$ws.Range("E46").Value = '  +5.66%  '
$ws.Range("E47").Value = '  -4.79%  '
$ws.Range("E48").Value = '  -1.68%  '
$ws.Range("D49").Value = '9.26'
$ws.Range("E49").Value = '  +5.69%  '
$ws.Range("D50").Value = '236.63'
$ws.Range("E50").Value = '  +15.16%  '
$ws.Range("D51").Value = '2.31'
$ws.Range("E51").Value = '  -2.42%  '

# Restore default style on the range (removes the temporary text-format override
# so the saved styles exactly match the original, unstyled Price/Volume cells).
$ws.Range("D2:E51").Style = "Normal"
